# Daily attendance processing - 2026-01-25 17:58:11
# Reorders the "Recorded By" (column G) comma-separated list of names/emails
# for the affected rows so that "System" remains last while the preceding
# entries are rotated: the leading entry is moved to sit immediately before
# "System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value2 = "system, backup@backdoor.com, System"
    }
}
